$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 7.813398806220115
$ws.Range("D2").Value = 16.22714636258608
$ws.Range("E2").Value = 6.462727774721766
$ws.Range("F2").Value = 118.6573573305936
$ws.Range("G2").Value = 4.137861705910042
$ws.Range("J2").Value = 11.72131370980549
$ws.Range("L2").Value = 8.733802619244567
$ws.Range("M2").Value = 62.73389755469431
$ws.Range("N2").Value = 18.21330573868288
$ws.Range("C3").Value = 7.850806573276888
$ws.Range("D3").Value = 16.13720891614603
$ws.Range("E3").Value = 6.236725123443867
$ws.Range("F3").Value = 118.6121778879305
$ws.Range("G3").Value = 4.15350076176412
$ws.Range("J3").Value = 11.79326505081509
$ws.Range("L3").Value = 8.712178421784669
$ws.Range("M3").Value = 62.30545170336337
$ws.Range("N3").Value = 18.23894920606617
$ws.Range("C4").Value = 7.874760783841378
$ws.Range("D4").Value = 16.08923392929637
$ws.Range("E4").Value = 6.093382676265841
$ws.Range("F4").Value = 118.6322417717166
$ws.Range("G4").Value = 4.163505240619097
$ws.Range("J4").Value = 11.83926718376379
$ws.Range("L4").Value = 8.701570291761398
$ws.Range("M4").Value = 62.06027000250345
$ws.Range("N4").Value = 18.2579213399345
$ws.Range("C5").Value = 7.884772535133008
$ws.Range("D5").Value = 16.07148745715095
$ws.Range("E5").Value = 6.033857302849843
$ws.Range("F5").Value = 118.6522452439535
$ws.Range("G5").Value = 4.167684547712224
$ws.Range("J5").Value = 11.85847660033368
$ws.Range("L5").Value = 8.697913601009562
$ws.Range("M5").Value = 61.96487781170207
$ws.Range("N5").Value = 18.26646883236016
$ws.Range("C6").Value = 7.886450159189902
$ws.Range("D6").Value = 16.06864901476864
$ws.Range("E6").Value = 6.02390720951613
$ws.Range("F6").Value = 118.6562754006829
$ws.Range("G6").Value = 4.168384740766714
$ws.Range("J6").Value = 11.86169441615701
$ws.Range("L6").Value = 8.697346493194516
$ws.Range("M6").Value = 61.94931166033444
$ws.Range("N6").Value = 18.26793756908074
$ws.Range("C7").Value = 7.874894789617088
$ws.Range("D7").Value = 16.08898731741058
$ws.Range("E7").Value = 6.092584344916893
$ws.Range("F7").Value = 118.6324639159835
$ws.Range("G7").Value = 4.163561187892806
$ws.Range("J7").Value = 11.8395243672894
$ws.Range("L7").Value = 8.701518285797473
$ws.Range("M7").Value = 62.05896517642626
$ws.Range("N7").Value = 18.25803330347265
$ws.Range("C8").Value = 7.826093862545941
$ws.Range("D8").Value = 16.19462142222376
$ws.Range("E8").Value = 6.385773339974516
$ws.Range("F8").Value = 118.6317425682623
$ws.Range("G8").Value = 4.143171408138285
$ws.Range("J8").Value = 11.74574681772001
$ws.Range("L8").Value = 8.725789722979513
$ws.Range("M8").Value = 62.58247330752025
$ws.Range("N8").Value = 18.22148075380806
$ws.Range("C9").Value = 7.738106364121719
$ws.Range("D9").Value = 16.45998475052729
$ws.Range("E9").Value = 6.922974484131649
$ws.Range("F9").Value = 119.0181116688683
$ws.Range("G9").Value = 4.106316018569705
$ws.Range("J9").Value = 11.57610559810305
$ws.Range("L9").Value = 8.79477235450363
$ws.Range("M9").Value = 63.74970125974507
$ws.Range("N9").Value = 18.17518976243209
$ws.Range("C10").Value = 7.678004418347897
$ws.Range("D10").Value = 16.69125949626301
$ws.Range("E10").Value = 7.293076024844286
$ws.Range("F10").Value = 119.5497476329953
$ws.Range("G10").Value = 4.081057752334118
$ws.Range("J10").Value = 11.45985315239418
$ws.Range("L10").Value = 8.858756290370611
$ws.Range("M10").Value = 64.69099592838778
$ws.Range("N10").Value = 18.15635566519108
$ws.Range("C11").Value = 7.651613051229248
$ws.Range("D11").Value = 16.8044874164912
$ws.Range("E11").Value = 7.455866207919808
$ws.Range("F11").Value = 119.8477532464
$ws.Range("G11").Value = 4.069942087340611
$ws.Range("J11").Value = 11.40871763225377
$ws.Range("L11").Value = 8.890802505142089
$ws.Range("M11").Value = 65.13689041335192
$ws.Range("N11").Value = 18.15101433815271
$ws.Range("C12").Value = 7.641752843053055
$ws.Range("D12").Value = 16.84852436674615
$ws.Range("E12").Value = 7.516694563011078
$ws.Range("F12").Value = 119.9688610611815
$ws.Range("G12").Value = 4.065785011212896
$ws.Range("J12").Value = 11.38959953694426
$ws.Range("L12").Value = 8.903363716698271
$ws.Range("M12").Value = 65.30823504183768
$ws.Range("N12").Value = 18.14944958826968
$ws.Range("C13").Value = 7.643870519144072
$ws.Range("D13").Value = 16.83898853968802
$ws.Range("E13").Value = 7.503630622544172
$ws.Range("F13").Value = 119.9424080951522
$ws.Range("G13").Value = 4.066678018237621
$ws.Range("J13").Value = 11.39370611932479
$ws.Range("L13").Value = 8.900639438244562
$ws.Range("M13").Value = 65.27122280455774
$ws.Range("N13").Value = 18.14976631931201
$ws.Range("C14").Value = 7.650799184721681
$ws.Range("D14").Value = 16.8080870818219
$ws.Range("E14").Value = 7.4608870384057
$ws.Range("F14").Value = 119.8575502865825
$ws.Range("G14").Value = 4.069599046097908
$ws.Range("J14").Value = 11.40713988127504
$ws.Range("L14").Value = 8.891827367721186
$ws.Range("M14").Value = 65.15093734844287
$ws.Range("N14").Value = 18.15087646393568
$ws.Range("C15").Value = 7.655060507623269
$ws.Range("D15").Value = 16.78931028522141
$ws.Range("E15").Value = 7.434598593881306
$ws.Range("F15").Value = 119.8066534537154
$ws.Range("G15").Value = 4.071395005482495
$ws.Range("J15").Value = 11.41540029150454
$ws.Range("L15").Value = 8.886485292449764
$ws.Range("M15").Value = 65.07758231133967
$ws.Range("N15").Value = 18.15161590694016
$ws.Range("C16").Value = 7.679748015132083
$ws.Range("D16").Value = 16.68402177970683
$ws.Range("E16").Value = 7.282323775582382
$ws.Range("F16").Value = 119.5314181584408
$ws.Range("G16").Value = 4.081791581794907
$ws.Range("J16").Value = 11.46322968631926
$ws.Range("L16").Value = 8.856721434193929
$ws.Range("M16").Value = 64.66220746784074
$ws.Range("N16").Value = 18.15676914474599
$ws.Range("C17").Value = 7.69513412778389
$ws.Range("D17").Value = 16.62149010419712
$ws.Range("E17").Value = 7.187468889658498
$ws.Range("F17").Value = 119.3770863771122
$ws.Range("G17").Value = 4.088264291472565
$ws.Range("J17").Value = 11.49301534242861
$ws.Range("L17").Value = 8.839217595192176
$ws.Range("M17").Value = 64.41188521748818
$ws.Range("N17").Value = 18.16075261586648
$ws.Range("C18").Value = 7.704073345061141
$ws.Range("D18").Value = 16.58627812070982
$ws.Range("E18").Value = 7.132387194911733
$ws.Range("F18").Value = 119.2935943049454
$ws.Range("G18").Value = 4.092022571243188
$ws.Range("J18").Value = 11.51031216877536
$ws.Range("L18").Value = 8.8294263598185
$ws.Range("M18").Value = 64.26957516616059
$ws.Range("N18").Value = 18.16334786448413
$ws.Range("C19").Value = 7.707115473633594
$ws.Range("D19").Value = 16.57448531571054
$ws.Range("E19").Value = 7.113648079185634
$ws.Range("F19").Value = 119.2662257398436
$ws.Range("G19").Value = 4.093301180692372
$ws.Range("J19").Value = 11.51619706637496
$ws.Range("L19").Value = 8.826158602240865
$ws.Range("M19").Value = 64.22167953417515
$ws.Range("N19").Value = 18.16427897439873
$ws.Range("C20").Value = 7.693487004091672
$ws.Range("D20").Value = 16.62806855740801
$ws.Range("E20").Value = 7.197620678136153
$ws.Range("F20").Value = 119.392967924737
$ws.Range("G20").Value = 4.087571613846607
$ws.Range("J20").Value = 11.48982758128447
$ws.Range("L20").Value = 8.841052266411644
$ws.Range("M20").Value = 64.43836010578949
$ws.Range("N20").Value = 18.16029714245816
$ws.Range("C21").Value = 7.648760464521624
$ws.Range("D21").Value = 16.8171320611814
$ws.Range("E21").Value = 7.47346413753741
$ws.Range("F21").Value = 119.8822494746941
$ws.Range("G21").Value = 4.068739666727905
$ws.Range("J21").Value = 11.40318743394379
$ws.Range("L21").Value = 8.894404099637937
$ws.Range("M21").Value = 65.18620082377393
$ws.Range("N21").Value = 18.15053800994264
$ws.Range("C22").Value = 7.620306563786228
$ws.Range("D22").Value = 16.947459421647
$ws.Range("E22").Value = 7.648978510648438
$ws.Range("F22").Value = 120.2502490522129
$ws.Range("G22").Value = 4.05673524916147
$ws.Range("J22").Value = 11.34799314786025
$ws.Range("L22").Value = 8.931757077423768
$ws.Range("M22").Value = 65.68947107313559
$ws.Range("N22").Value = 18.14682582180229
$ws.Range("C23").Value = 7.635422787846091
$ws.Range("D23").Value = 16.87728056232197
$ws.Range("E23").Value = 7.555743404745489
$ws.Range("F23").Value = 120.0493683307539
$ws.Range("G23").Value = 4.063115044703283
$ws.Range("J23").Value = 11.3773224285009
$ws.Range("L23").Value = 8.911592727377315
$ws.Range("M23").Value = 65.41955531120242
$ws.Range("N23").Value = 18.14856521726656
$ws.Range("C24").Value = 7.694231377151907
$ws.Range("D24").Value = 16.6250921395678
$ws.Range("E24").Value = 7.193032761343029
$ws.Range("F24").Value = 119.3857715841263
$ws.Range("G24").Value = 4.087884657927585
$ws.Range("J24").Value = 11.4912682298933
$ws.Range("L24").Value = 8.840221964726098
$ws.Range("M24").Value = 64.42638581028967
$ws.Range("N24").Value = 18.16050211121814
$ws.Range("C25").Value = 7.761099862064076
$ws.Range("D25").Value = 16.3818703004049
$ws.Range("E25").Value = 6.781867606926054
$ws.Range("F25").Value = 118.8708823009834
$ws.Range("G25").Value = 4.115960307278533
$ws.Range("J25").Value = 11.62050307005564
$ws.Range("L25").Value = 8.773792320803455
$ws.Range("M25").Value = 63.4190679464839
$ws.Range("N25").Value = 18.18502505134802
